$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arr = New-Object 'object[,]' 19,6
$arr[0,0] = 1.455362044514542
$arr[0,1] = 0.306821227259698
$arr[0,2] = 0.1494219747398047
$arr[0,3] = 0.4942365360607697
$arr[0,4] = 1
$arr[0,5] = 2.405841782574814
$arr[1,0] = 1.455362044514542
$arr[1,1] = 1.655778082260271
$arr[1,2] = 0.7527432677738641
$arr[1,3] = 0.4942365360607697
$arr[1,4] = 0
$arr[1,5] = 4.358119930609447
$arr[2,0] = 3.286832544864788
$arr[2,1] = 1.655778082260271
$arr[2,2] = 3.537761648806719
$arr[2,3] = 0.4942365360607697
$arr[2,4] = 1
$arr[2,5] = 8.974608811992548
$arr[3,0] = 0.2917716402565462
$arr[3,1] = 1.655778082260271
$arr[3,2] = 0.1494219747398047
$arr[3,3] = 0.4942365360607697
$arr[3,4] = 1
$arr[3,5] = 2.591208233317391
$arr[4,0] = 3.286832544864788
$arr[4,1] = 1.655778082260271
$arr[4,2] = 0.7527432677738641
$arr[4,3] = 0.4942365360607697
$arr[4,4] = 1
$arr[4,5] = 6.189590430959694
$arr[5,0] = 3.286832544864788
$arr[5,1] = 1.655778082260271
$arr[5,2] = 0.1494219747398047
$arr[5,3] = 0.4942365360607697
$arr[5,4] = 1
$arr[5,5] = 5.586269137925634
$arr[6,0] = 3.286832544864788
$arr[6,1] = 1.655778082260271
$arr[6,2] = 3.537761648806719
$arr[6,3] = 0.4942365360607697
$arr[6,4] = 1
$arr[6,5] = 8.974608811992548
$arr[7,0] = 3.286832544864788
$arr[7,1] = 1.655778082260271
$arr[7,2] = 0.7527432677738641
$arr[7,3] = 0.4942365360607697
$arr[7,4] = 1
$arr[7,5] = 6.189590430959694
$arr[8,0] = 3.286832544864788
$arr[8,1] = 1.655778082260271
$arr[8,2] = 0.7527432677738641
$arr[8,3] = 0.4942365360607697
$arr[8,4] = 1
$arr[8,5] = 6.189590430959694
$arr[9,0] = 3.286832544864788
$arr[9,1] = 1.655778082260271
$arr[9,2] = 0.7527432677738641
$arr[9,3] = 0.4942365360607697
$arr[9,4] = 1
$arr[9,5] = 6.189590430959694
$arr[10,0] = 0.1190320826869504
$arr[10,1] = 0.306821227259698
$arr[10,2] = 0.7527432677738641
$arr[10,3] = 0.4942365360607697
$arr[10,4] = 1
$arr[10,5] = 1.672833113781282
$arr[11,0] = 0.04271373187048222
$arr[11,1] = 0.306821227259698
$arr[11,2] = 0.1494219747398047
$arr[11,3] = 0.4942365360607697
$arr[11,4] = 1
$arr[11,5] = 0.9931934699307545
$arr[12,0] = 3.286832544864788
$arr[12,1] = 117.745847958593
$arr[12,2] = 22.3905356188092
$arr[12,3] = 10.19245300693656
$arr[12,4] = 0
$arr[12,5] = 153.6156691292036
$arr[13,0] = 3.286832544864788
$arr[13,1] = 1.655778082260271
$arr[13,2] = 3.537761648806719
$arr[13,3] = 0.4942365360607697
$arr[13,4] = 1
$arr[13,5] = 8.974608811992548
$arr[14,0] = 3.286832544864788
$arr[14,1] = 1.655778082260271
$arr[14,2] = 0.7527432677738641
$arr[14,3] = 0.4942365360607697
$arr[14,4] = 1
$arr[14,5] = 6.189590430959694
$arr[15,0] = 1.455362044514542
$arr[15,1] = 1.655778082260271
$arr[15,2] = 0.1494219747398047
$arr[15,3] = 0.4942365360607697
$arr[15,4] = 1
$arr[15,5] = 3.754798637575387
$arr[16,0] = 1.455362044514542
$arr[16,1] = 1.655778082260271
$arr[16,2] = 0.1494219747398047
$arr[16,3] = 0.4942365360607697
$arr[16,4] = 1
$arr[16,5] = 3.754798637575387
$arr[17,0] = 0.6606524410359556
$arr[17,1] = 1.655778082260271
$arr[17,2] = 3.537761648806719
$arr[17,3] = 0.4942365360607697
$arr[17,4] = 0
$arr[17,5] = 6.348428708163715
$arr[18,0] = 0.0006408296065709695
$arr[18,1] = 0.002571899574220771
$arr[18,2] = 0.7527432677738641
$arr[18,3] = 0.4942365360607697
$arr[18,4] = 0
$arr[18,5] = 1.250192533015426

$ws.Range("B2:G20").Value = $arr
